# All final evaluation results
# Update Summary, Classification Report, and Confusion Matrix sheets with the
# latest run's metrics for isolation_forest/augmented/none_3/split_4/test_95_5.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.0498220640569395
$summary.Range("C2").Value = 0.0498220640569395
$summary.Range("D2").Value = 1
$summary.Range("E2").Value = 0.09491525423728814
$summary.Range("F2").Value = 0.2077151335311573
$summary.Range("G2").Value = 0.5768621236133122
$summary.Range("H2").Value = 0.7486623863028358
$summary.Range("I2").Value = 28
$summary.Range("J2").Value = 534
$summary.Range("K2").Value = 0
$summary.Range("L2").Value = 0

# --- Classification Report sheet ---
$report = $wb.Worksheets.Item("Classification Report")

# label "0"
$report.Range("B2").Value = 0
$report.Range("C2").Value = 0
$report.Range("D2").Value = 0

# label "1"
$report.Range("B3").Value = 0.0498220640569395
$report.Range("C3").Value = 1
$report.Range("D3").Value = 0.09491525423728814

# accuracy
$report.Range("B4").Value = 0.0498220640569395
$report.Range("C4").Value = 0.0498220640569395
$report.Range("D4").Value = 0.0498220640569395
$report.Range("E4").Value = 0.0498220640569395

# macro avg
$report.Range("B5").Value = 0.02491103202846975
$report.Range("C5").Value = 0.5
$report.Range("D5").Value = 0.04745762711864407

# weighted avg
$report.Range("B6").Value = 0.002482238066893783
$report.Range("C6").Value = 0.0498220640569395
$report.Range("D6").Value = 0.004728873876590867

# --- Confusion Matrix sheet ---
$confusion = $wb.Worksheets.Item("Confusion Matrix")

# Actual 0
$confusion.Range("B2").Value = 0
$confusion.Range("C2").Value = 534

# Actual 1
$confusion.Range("B3").Value = 0
$confusion.Range("C3").Value = 28
